$d = $word.ActiveDocument

# --- Locate the "Malaysian Penal Code" paragraph (start of the block to clear) ---
$findRange = $d.Content
$found = $findRange.Find.Execute("Malaysian Penal Code")
if (-not $found) {
    throw "Could not find 'Malaysian Penal Code' text"
}
$targetStart = $findRange.Start

$count = $d.Paragraphs.Count
$startIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $targetStart -and $p.Range.End -gt $targetStart) {
        $startIdx = $i
        break
    }
}
if ($startIdx -eq -1) {
    throw "Could not resolve paragraph index for 'Malaysian Penal Code'"
}

# The five paragraphs involved (by their original layout):
#   startIdx+0 : "Malaysian Penal Code"
#   startIdx+1 : (empty)
#   startIdx+2 : "Penal Code is a law that "
#   startIdx+3 : (empty)
#   startIdx+4 : "According to the law" + " to Section 147" + ", such an offence...crime"
# They must become SIX empty paragraphs (same <w:pPr> spacing, no runs).

function Clear-ParaText($idx) {
    $p = $d.Paragraphs.Item($idx)
    $r = $p.Range
    if ($r.End - 1 -gt $r.Start) {
        $clearRange = $d.Range($r.Start, $r.End - 1)
        $clearRange.Delete()
    }
}

# Split the first paragraph's text in two (while text is still present) so the
# split produces a clean extra paragraph mark with no leftover empty run -
# this yields the 6th paragraph that the edit introduces.
$p0 = $d.Paragraphs.Item($startIdx)
$p0Start = $p0.Range.Start
$splitPoint = $d.Range($p0Start + 9, $p0Start + 9)
$splitPoint.InsertParagraphAfter()

# After the split the paragraphs shift by one:
#   startIdx+0 : "Malaysian"
#   startIdx+1 : " Penal Code"
#   startIdx+2 : (empty, was startIdx+1)
#   startIdx+3 : "Penal Code is a law that " (was startIdx+2)
#   startIdx+4 : (empty, was startIdx+3)
#   startIdx+5 : "According to the law..." (was startIdx+4)

# Clear remaining text, from the last affected paragraph back to the first so
# that earlier indices stay valid while later ones are processed.
Clear-ParaText ($startIdx + 5)
Clear-ParaText ($startIdx + 3)
Clear-ParaText ($startIdx + 1)
Clear-ParaText ($startIdx + 0)
